$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simulation score")

$ws.Range("Q2").Value = 13360.52267
$ws.Range("R2").Value = 13360.52267
$ws.Range("S2").Value = 13821
$ws.Range("E3").Value = 27
$ws.Range("F3").Value = 2925.6
$ws.Range("G3").Value = 3
$ws.Range("Q3").Value = 11164.4
$ws.Range("R3").Value = 11164.4
$ws.Range("S3").Value = 13032.875330000001
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = 2656.6
$ws.Range("G4").Value = 3
$ws.Range("Q4").Value = 13198.17489
$ws.Range("R4").Value = 13198.17489
$ws.Range("S4").Value = 11164.4
$ws.Range("Q5").Value = 13821
$ws.Range("R5").Value = 13821
$ws.Range("S5").Value = 11164.4
$ws.Range("Q6").Value = 11593.883669999999
$ws.Range("R6").Value = 11593.883669999999
$ws.Range("S6").Value = 11893.48911
$ws.Range("Q7").Value = 12077.97522
$ws.Range("R7").Value = 12077.97522
$ws.Range("S7").Value = 13821
$ws.Range("AC8").Value = 0
$ws.Range("AD8").Value = 0
$ws.Range("Q8").Value = 13821
$ws.Range("R8").Value = 13821
$ws.Range("S8").Value = 13821
$ws.Range("Q9").Value = 13375.281559999999
$ws.Range("R9").Value = 13375.281559999999
$ws.Range("S9").Value = 13821
$ws.Range("Q10").Value = 11164.4
$ws.Range("R10").Value = 11164.4
$ws.Range("S10").Value = 13772.29567
$ws.Range("Q11").Value = 13115.52511
$ws.Range("R11").Value = 13115.52511
$ws.Range("S11").Value = 11164.4
$ws.Range("Q12").Value = 13821
$ws.Range("R12").Value = 13821
$ws.Range("S12").Value = 11164.4
$ws.Range("Q13").Value = 12126.67956
$ws.Range("R13").Value = 12126.67956
$ws.Range("S13").Value = 11608.64256
$ws.Range("Q14").Value = 12714.083329999999
$ws.Range("R14").Value = 12714.083329999999
$ws.Range("S14").Value = 13821
$ws.Range("Q15").Value = 11841.343999999999
$ws.Range("R15").Value = 11841.343999999999
$ws.Range("S15").Value = 10781.62667
$ws.Range("S16").Value = 4.4276666669999996
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("Q21").Value = 670.05355559999998
$ws.Range("R21").Value = 670.05355559999998
$ws.Range("S21").Value = 875.20211110000002
$ws.Range("Q22").Value = 13821
$ws.Range("R22").Value = 13821
$ws.Range("S22").Value = 13817.89878
$ws.Range("Q23").Value = 12875.056
$ws.Range("R23").Value = 12875.056
$ws.Range("S23").Value = 12668.63867
$ws.Range("Q24").Value = 10557.685219999999
$ws.Range("R24").Value = 10557.685219999999
$ws.Range("S24").Value = 10895.4
$ws.Range("Q25").Value = 11050.62667
$ws.Range("R25").Value = 11050.62667
$ws.Range("S25").Value = 9709.6564440000002
$ws.Range("Q26").Value = 13351.66733
$ws.Range("R26").Value = 13351.66733
$ws.Range("S26").Value = 11164.4
$ws.Range("Q27").Value = 13821
$ws.Range("R27").Value = 13821
$ws.Range("S27").Value = 11164.4
$ws.Range("Q28").Value = 13147.99467
$ws.Range("R28").Value = 13147.99467
$ws.Range("S28").Value = 12659.47544
$ws.Range("Q29").Value = 11164.4
$ws.Range("R29").Value = 11164.4
$ws.Range("S29").Value = 13821
$ws.Range("Q30").Value = 12753.93233
$ws.Range("R30").Value = 12753.93233
$ws.Range("S30").Value = 13821
$ws.Range("Q31").Value = 13821
$ws.Range("R31").Value = 13821
$ws.Range("S31").Value = 13821
$ws.Range("Q32").Value = 12786.401889999999
$ws.Range("R32").Value = 12786.401889999999
$ws.Range("S32").Value = 12184.239219999999
$ws.Range("Q33").Value = 11164.4
$ws.Range("R33").Value = 11164.4
$ws.Range("S33").Value = 11164.4
$ws.Range("Q34").Value = 13748.68144
$ws.Range("R34").Value = 13748.68144
$ws.Range("S34").Value = 11164.4
$ws.Range("Q35").Value = 13821
$ws.Range("R35").Value = 13821
$ws.Range("S35").Value = 13211.45789
$ws.Range("Q36").Value = 8640.6579999999994
$ws.Range("R36").Value = 8640.6579999999994
$ws.Range("S36").Value = 12135.529329999999
$ws.Range("Q37").Value = 1908.324333
$ws.Range("R37").Value = 1908.324333
$ws.Range("S37").Value = 4809.9218890000002
$ws.Range("S38").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = 0
$ws.Range("Q42").Value = 13821
$ws.Range("R42").Value = 13821
$ws.Range("S42").Value = 12368.57589
$ws.Range("Q43").Value = 11875.77844
$ws.Range("R43").Value = 11875.77844
$ws.Range("S43").Value = 13821
$ws.Range("Q44").Value = 13339.86022
$ws.Range("R44").Value = 13339.86022
$ws.Range("S44").Value = 13821
$ws.Range("Q45").Value = 13821
$ws.Range("R45").Value = 13821
$ws.Range("S45").Value = 13821
$ws.Range("Q46").Value = 13821
$ws.Range("R46").Value = 13821
$ws.Range("S46").Value = 13821
$ws.Range("Q47").Value = 12777.546560000001
$ws.Range("R47").Value = 12777.546560000001
$ws.Range("S47").Value = 13392.99222
$ws.Range("Q48").Value = 11356.26556
$ws.Range("R48").Value = 11356.26556
$ws.Range("S48").Value = 11164.4
$ws.Range("Q49").Value = 13821
$ws.Range("R49").Value = 13821
$ws.Range("S49").Value = 11164.4
